$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intInit")

# Rows 2-8 keep their B values; only the label in column A shifts up by one
# slot because the "allianceFightTimePerFight" row was removed from the
# underlying lookup table. Row 6 also gets a new B value (7 instead of 60).
$ws.Range("A2").Value = "createAllianceGem"
$ws.Range("A3").Value = "buyArchonGem"
$ws.Range("A4").Value = "editAllianceBasicInfoGem"
$ws.Range("A5").Value = "editAllianceTerrianHonour"
$ws.Range("A6").Value = "activeShrineStageEventTime"
$ws.Range("B6").Value = 7
$ws.Range("A7").Value = "allianceFightPrepareTime"

# Row 8 (allianceFightTotalFightTime / 300) is unchanged.

# Rows 9-12 are relabeled/renumbered and row 9's value changes.
$ws.Range("A9").Value = "allianceHelpDefenceTroopsMaxCount"
$ws.Range("B9").Value = 2
$ws.Range("A10").Value = "allianceRevengeMaxTime"
$ws.Range("B10").Value = 300
$ws.Range("A11").Value = "dragonStrikeHpDecreasedPercent"
$ws.Range("B11").Value = 5
$ws.Range("A12").Value = "allianceRegionMapBaseTimePerGrid"
$ws.Range("B12").Value = 1

# Two brand-new rows for the alliance region map dimensions.
$ws.Rows.Item(13).RowHeight = 20
$ws.Range("A13").Value = "allianceRegionMapWidth"
$ws.Range("B13").Value = 51

$ws.Rows.Item(14).RowHeight = 20
$ws.Range("A14").Value = "allianceRegionMapHeight"
$ws.Range("B14").Value = 51
